# Generate Report for Handoff
#
# A new handoff was generated: the localized file's GUID-based name changed
# (new content -> new hash), and the handoff/target timestamps were
# refreshed to reflect the new run.

$wb = $excel.ActiveWorkbook

# old guid: 34fbaffe-da07-4da8-8091-71b7ec256a63
# old content hash: 265c527b48ba0c8fb30de6fd61caaddf90438cd2
$newGuid = "9d7dd74f-3ea4-4102-aea8-78ae5ee23819"
$newHash = "4aadf6bcab74acbe0225699c6b4178e3730c5f87"

$newFileName = "$newGuid.md"
$newPathName = "e2e\$newGuid.md"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathName
$wsOverview.Range("G2").Value = "2016-09-06 11:17:37"
foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = $newPathName
}

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newFileName
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-06 11:17:32"
foreach ($h in $wsZhCn.Hyperlinks) {
    $h.TextToDisplay = $newFileName
}

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newFileName
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-06 11:17:37"
foreach ($h in $wsDeDe.Hyperlinks) {
    $h.TextToDisplay = $newFileName
}
